# Carbon_Sediment_Samples.xlsx edit
# "viendo carbon concentrations for jesus"
#
# Summary of the change (from the OOXML diff):
#  - On the "All Samples" worksheet, insert a brand-new blank column in front
#    of column N. Everything that used to live in N..S shifts right to O..T
#    (this also shifts the chart series formulas and the floating chart
#    anchors on that sheet, which Excel does automatically on a column
#    insert).
#  - Populate the freshly vacated M:N cells (rows 2-9) with a small
#    avg/min/max mgC summary block that reads off the existing E and K
#    columns.
#  - A handful of view-state tweaks (zoom, scroll position, selection,
#    active sheet) on the four sheets.

$wb = $excel.ActiveWorkbook

$wsAll = $wb.Worksheets.Item("All Samples")
$wsComposite = $wb.Worksheets.Item("Composite Samples")
$wsSeasonal = $wb.Worksheets.Item("Seasonal")
$wsSignificance = $wb.Worksheets.Item("Significance")

# ---------------------------------------------------------------------
# 1. Insert a new column before column N (14) on "All Samples".
#    This shifts the existing N:S mini-tables to O:T and drags every
#    chart series / drawing anchor on the sheet along with it.
# ---------------------------------------------------------------------
$wsAll.Activate()
$wsAll.Columns.Item(14).Insert()

# ---------------------------------------------------------------------
# 2. New avg / min / max mgC summary block in M2:N9.
# ---------------------------------------------------------------------
$wsAll.Range("M2").HorizontalAlignment = -4108
$wsAll.Range("M2").Value = "avg mgC"

$wsAll.Range("N2").HorizontalAlignment = -4108
$wsAll.Range("N2").Value = "avg mgC/mg"

$wsAll.Range("M3").NumberFormat = "0.0000"
$wsAll.Range("M3").HorizontalAlignment = -4108
$wsAll.Range("M3").Formula = "=AVERAGE(E2:E70)"

$wsAll.Range("N3").NumberFormat = "0.0000"
$wsAll.Range("N3").HorizontalAlignment = -4108
$wsAll.Range("N3").Formula = "=AVERAGE(K2:K70)"

$wsAll.Range("M5").HorizontalAlignment = -4108
$wsAll.Range("M5").Value = "min mgC"

$wsAll.Range("N5").HorizontalAlignment = -4108
$wsAll.Range("N5").Value = "min mgC/mg"

$wsAll.Range("M6").NumberFormat = "0.0000"
$wsAll.Range("M6").HorizontalAlignment = -4108
$wsAll.Range("M6").Formula = "=MIN(E2:E70)"

$wsAll.Range("N6").NumberFormat = "0.0000"
$wsAll.Range("N6").HorizontalAlignment = -4108
$wsAll.Range("N6").Formula = "=MIN(K2:K70)"

$wsAll.Range("M8").HorizontalAlignment = -4108
$wsAll.Range("M8").Value = "max mgC"

$wsAll.Range("N8").HorizontalAlignment = -4108
$wsAll.Range("N8").Value = "max mgC/mg"

$wsAll.Range("M9").NumberFormat = "0.0000"
$wsAll.Range("M9").HorizontalAlignment = -4108
$wsAll.Range("M9").Formula = "=MAX(E2:E70)"

$wsAll.Range("N9").NumberFormat = "0.0000"
$wsAll.Range("N9").HorizontalAlignment = -4108
$wsAll.Range("N9").Formula = "=MAX(K2:K70)"

# ---------------------------------------------------------------------
# 3. View-state housekeeping matching where the author was last looking.
#    (topLeftCell/scroll position is best-effort; selection + zoom +
#    active sheet are what actually round-trip.)
# ---------------------------------------------------------------------
$wsComposite.Activate()
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 10

$wsSeasonal.Activate()
$excel.ActiveWindow.ScrollRow = 52
$excel.ActiveWindow.ScrollColumn = 1

$wsSignificance.Activate()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1

$wsAll.Activate()
$excel.ActiveWindow.Zoom = 72
$excel.ActiveWindow.ScrollRow = 44
$excel.ActiveWindow.ScrollColumn = 1
$wsAll.Range("K55").Select()
